# Update the runs/balls/fours/sixes figures for Ravindra Jadeja's Chennai
# Super Kings rows. The source sheet stores these numeric-looking figures
# as TEXT (not numbers), so each cell is written as text rather than a
# number. We temporarily format the target range as Text ("@") so the
# COM layer does not auto-coerce the literal into a numeric cell, then
# restore the range's style back to "Normal" afterwards so no residual
# number-format/style change is left behind - only the cell values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @("12", "9", "1", "0")
    5  = @("50", "35", "5", "2")
    6  = @("25", "10", "3", "1")
    7  = @("33", "13", "0", "4")
    8  = @("31", "11", "2", "3")
    9  = @("10", "5", "2", "0")
    10 = @("1", "2", "0", "0")
    12 = @("35", "30", "4", "0")
}

foreach ($row in $updates.Keys) {
    $values = $updates[$row]
    $rowRange = $ws.Range("C$row`:F$row")

    # Force text storage for this row's C:F cells.
    $rowRange.NumberFormat = "@"

    $ws.Range("C$row").Value = $values[0]
    $ws.Range("D$row").Value = $values[1]
    $ws.Range("E$row").Value = $values[2]
    $ws.Range("F$row").Value = $values[3]

    # Drop the temporary format back to Normal so no style diff remains.
    $rowRange.Style = "Normal"
}
